$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency market data (price / 1h volume change,
# plus two swapped coin rows) as produced by the scheduled GitHub Actions run.
# Values are written with a leading apostrophe to force Excel to keep them
# as literal text (matching the original inlineStr/text cells) instead of
# reinterpreting number-like strings (e.g. "1.00", "46.016.50") as numbers.

$ws.Range("D2").Value = "'46.016.50"
$ws.Range("E2").Value = "'  -0.80%  "
$ws.Range("D3").Value = "'2.596.57"
$ws.Range("E3").Value = "'  +0.63%  "
$ws.Range("E4").Value = "'  +0.16%  "
$ws.Range("D5").Value = "'312.02"
$ws.Range("E5").Value = "'  +2.47%  "
$ws.Range("D6").Value = "'98.89"
$ws.Range("E6").Value = "'  -0.67%  "
$ws.Range("D7").Value = "'0.596"
$ws.Range("E7").Value = "'  -0.44%  "
$ws.Range("E8").Value = "'  +0.18%  "
$ws.Range("D9").Value = "'0.578"
$ws.Range("E9").Value = "'  +0.15%  "
$ws.Range("D10").Value = "'38.86"
$ws.Range("E10").Value = "'  +0.93%  "
$ws.Range("D11").Value = "'54.41"
$ws.Range("E11").Value = "'  -0.86%  "
$ws.Range("D12").Value = "'0.0837"
$ws.Range("E12").Value = "'  +0.25%  "
$ws.Range("D13").Value = "'8.11"
$ws.Range("E13").Value = "'  -0.30%  "
$ws.Range("D14").Value = "'3.004.74"
$ws.Range("E14").Value = "'  +0.92%  "
$ws.Range("E15").Value = "'  +1.05%  "
$ws.Range("D16").Value = "'2.616.55"
$ws.Range("E16").Value = "'  +0.67%  "
$ws.Range("D17").Value = "'0.914"
$ws.Range("E17").Value = "'  +1.53%  "
$ws.Range("D18").Value = "'14.81"
$ws.Range("E18").Value = "'  -0.15%  "
$ws.Range("D19").Value = "'46.193.94"
$ws.Range("E19").Value = "'  -0.61%  "
$ws.Range("E20").Value = "'  +0.63%  "
$ws.Range("D21").Value = "'12.81"
$ws.Range("E21").Value = "'  -3.82%  "
$ws.Range("D22").Value = "'6.71"
$ws.Range("E22").Value = "'  +1.18%  "
$ws.Range("B23").Value = "'BitcoinCash"
$ws.Range("C23").Value = "'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").Value = "'276.65"
$ws.Range("E23").Value = "'  +8.51%  "
$ws.Range("B24").Value = "'Litecoin"
$ws.Range("C24").Value = "'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "'71.80"
$ws.Range("E24").Value = "'  +0.94%  "
$ws.Range("E25").Value = "'  +3.44%  "
$ws.Range("E26").Value = "'  +1.72%  "
$ws.Range("D27").Value = "'29.95"
$ws.Range("E27").Value = "'  +7.87%  "
$ws.Range("E28").Value = "'  +0.00%  "
$ws.Range("E29").Value = "'  +1.35%  "
$ws.Range("D30").Value = "'10.71"
$ws.Range("E30").Value = "'  +2.49%  "
$ws.Range("D31").Value = "'2.20"
$ws.Range("E31").Value = "'  -4.07%  "
$ws.Range("D32").Value = "'37.80"
$ws.Range("E32").Value = "'  -4.47%  "
$ws.Range("D33").Value = "'6.22"
$ws.Range("E33").Value = "'  +1.79%  "
$ws.Range("D34").Value = "'3.59"
$ws.Range("E34").Value = "'  -3.93%  "
$ws.Range("D35").Value = "'155.84"
$ws.Range("E35").Value = "'  +4.04%  "
$ws.Range("D36").Value = "'2.21"
$ws.Range("E36").Value = "'  -4.42%  "
$ws.Range("D37").Value = "'0.0835"
$ws.Range("E37").Value = "'  +0.29%  "
$ws.Range("E38").Value = "'  -4.37%  "
$ws.Range("E39").Value = "'  +5.78%  "
$ws.Range("E40").Value = "'  +1.00%  "
$ws.Range("D41").Value = "'23.51"
$ws.Range("E41").Value = "'  +29.55%  "
$ws.Range("D42").Value = "'15.73"
$ws.Range("E42").Value = "'  -0.03%  "
$ws.Range("D43").Value = "'0.0330"
$ws.Range("E43").Value = "'  +2.40%  "
$ws.Range("D44").Value = "'3.58"
$ws.Range("E44").Value = "'  -1.27%  "
$ws.Range("D45").Value = "'3.96"
$ws.Range("E45").Value = "'  -5.02%  "
$ws.Range("D46").Value = "'2.096.20"
$ws.Range("E46").Value = "'  +3.31%  "
$ws.Range("D47").Value = "'1.00"
$ws.Range("E47").Value = "'  +0.09%  "
$ws.Range("D48").Value = "'95.02"
$ws.Range("E48").Value = "'  +3.98%  "
$ws.Range("D49").Value = "'9.58"
$ws.Range("E49").Value = "'  +6.03%  "
$ws.Range("D50").Value = "'108.47"
$ws.Range("E50").Value = "'  +0.00%  "
$ws.Range("B51").Value = "'RocketPoolETH"
$ws.Range("C51").Value = "'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "'2.857.50"
$ws.Range("E51").Value = "'  +0.70%  "
